$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells (coin name / link) - rows 48 and 49 swapped values
$plainUpdates = @(
    @{ Cell = "B48"; Value = "TrustWalletToken" },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt" },
    @{ Cell = "B49"; Value = "BinanceUSD" },
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd" }
)

foreach ($u in $plainUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Price (D) and Volume(1h) (E) cells - stored as text, force Text format
# so Excel does not coerce numeric-looking strings into numbers/floats,
# then reset the style back to Normal so no stray style index is left behind.
$textUpdates = @(
    @{ Cell = "D2"; Value = "43.738.52" },
    @{ Cell = "E2"; Value = "  -0.09%  " },
    @{ Cell = "D3"; Value = "2.351.82" },
    @{ Cell = "E3"; Value = "  +0.55%  " },
    @{ Cell = "E4"; Value = "  -0.22%  " },
    @{ Cell = "D5"; Value = "239.64" },
    @{ Cell = "E5"; Value = "  +0.83%  " },
    @{ Cell = "D6"; Value = "0.669" },
    @{ Cell = "E6"; Value = "  +0.04%  " },
    @{ Cell = "D7"; Value = "74.40" },
    @{ Cell = "E7"; Value = "  +2.30%  " },
    @{ Cell = "E8"; Value = "  -0.05%  " },
    @{ Cell = "E9"; Value = "  +2.68%  " },
    @{ Cell = "E10"; Value = "  +3.47%  " },
    @{ Cell = "D11"; Value = "59.91" },
    @{ Cell = "E11"; Value = "  +4.82%  " },
    @{ Cell = "D12"; Value = "37.26" },
    @{ Cell = "E12"; Value = "  +16.59%  " },
    @{ Cell = "D13"; Value = "7.35" },
    @{ Cell = "E13"; Value = "  +2.57%  " },
    @{ Cell = "E14"; Value = "  +1.11%  " },
    @{ Cell = "D15"; Value = "2.703.11" },
    @{ Cell = "E15"; Value = "  +0.66%  " },
    @{ Cell = "D16"; Value = "16.36" },
    @{ Cell = "E16"; Value = "  +0.44%  " },
    @{ Cell = "E17"; Value = "  +4.38%  " },
    @{ Cell = "D18"; Value = "2.352.13" },
    @{ Cell = "E18"; Value = "  +0.89%  " },
    @{ Cell = "D19"; Value = "43.720.72" },
    @{ Cell = "E19"; Value = "  +0.30%  " },
    @{ Cell = "D20"; Value = "0.0000104" },
    @{ Cell = "E20"; Value = "  +3.47%  " },
    @{ Cell = "D21"; Value = "6.60" },
    @{ Cell = "E21"; Value = "  -3.03%  " },
    @{ Cell = "D22"; Value = "77.13" },
    @{ Cell = "E22"; Value = "  +0.97%  " },
    @{ Cell = "D23"; Value = "253.26" },
    @{ Cell = "E23"; Value = "  -0.67%  " },
    @{ Cell = "E24"; Value = "  +0.08%  " },
    @{ Cell = "D25"; Value = "3.79" },
    @{ Cell = "E25"; Value = "  +3.74%  " },
    @{ Cell = "E26"; Value = "  -4.50%  " },
    @{ Cell = "E27"; Value = "  +1.21%  " },
    @{ Cell = "D28"; Value = "10.69" },
    @{ Cell = "E28"; Value = "  +1.59%  " },
    @{ Cell = "D29"; Value = "2.30" },
    @{ Cell = "E29"; Value = "  +0.68%  " },
    @{ Cell = "D30"; Value = "22.28" },
    @{ Cell = "E30"; Value = "  -0.77%  " },
    @{ Cell = "D31"; Value = "174.82" },
    @{ Cell = "E31"; Value = "  +0.32%  " },
    @{ Cell = "E32"; Value = "  +0.70%  " },
    @{ Cell = "E33"; Value = "  -0.57%  " },
    @{ Cell = "D34"; Value = "0.0756" },
    @{ Cell = "E34"; Value = "  +1.32%  " },
    @{ Cell = "D35"; Value = "5.53" },
    @{ Cell = "E35"; Value = "  +0.02%  " },
    @{ Cell = "D36"; Value = "5.14" },
    @{ Cell = "E36"; Value = "  +0.30%  " },
    @{ Cell = "D37"; Value = "3.79" },
    @{ Cell = "E37"; Value = "  +2.34%  " },
    @{ Cell = "D38"; Value = "6.59" },
    @{ Cell = "E38"; Value = "  +6.73%  " },
    @{ Cell = "E39"; Value = "  +2.08%  " },
    @{ Cell = "E40"; Value = "  +1.86%  " },
    @{ Cell = "D41"; Value = "5.53" },
    @{ Cell = "E41"; Value = "  +18.23%  " },
    @{ Cell = "D42"; Value = "21.18" },
    @{ Cell = "E42"; Value = "  +14.26%  " },
    @{ Cell = "D43"; Value = "65.97" },
    @{ Cell = "E43"; Value = "  +10.56%  " },
    @{ Cell = "E44"; Value = "  -1.86%  " },
    @{ Cell = "D45"; Value = "9.05" },
    @{ Cell = "E45"; Value = "  +1.92%  " },
    @{ Cell = "D46"; Value = "0.202" },
    @{ Cell = "E46"; Value = "  +0.35%  " },
    @{ Cell = "E47"; Value = "  +3.62%  " },
    @{ Cell = "D48"; Value = "1.25" },
    @{ Cell = "E48"; Value = "  +1.92%  " },
    @{ Cell = "D49"; Value = "1.00" },
    @{ Cell = "E49"; Value = "  +0.00%  " },
    @{ Cell = "D50"; Value = "1.16" },
    @{ Cell = "E50"; Value = "  +1.35%  " },
    @{ Cell = "D51"; Value = "98.00" }
)

foreach ($u in $textUpdates) {
    $r = $ws.Range($u.Cell)
    $r.NumberFormat = "@"
    $r.Value = $u.Value
    $r.Style = "Normal"
}
